$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")

# --- Insert two new rows right after row 45 (before the old row 46) -------
# Row 45 currently holds a placeholder "0" entry under the 2023-10-10 date
# header (row 44). The edit turns that placeholder into a real entry and
# adds two more entries for the same date, pushing every following row
# down by two. Inserting above row 46 shifts all subsequent rows/merges
# and keeps the trailing date sequence identical, just relocated two rows
# further down (45210, 45211, ... end up two rows lower than before).
[void]$ws.Rows("46:47").Insert()

# Make sure the two new rows look like the other "entry" rows (same number
# format / fill / border / alignment as row 45, an existing entry row
# under the same date header) instead of whatever Insert() guessed.
[void]$ws.Range("A45:B45").Copy()
[void]$ws.Range("A46:B47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the three entries for the 2023-10-10 (A44) date block -------
$ws.Range("A45").Value2 = 10
$ws.Range("B45").Value2 = "DOM ORIONE"

$ws.Range("A46").Value2 = 15
$ws.Range("B46").Value2 = "VILA AZUL"

$ws.Range("A47").Value2 = 15
$ws.Range("B47").Value2 = "LAGO AZUL"

# --- Fix up the TOTAL formula/row (now row 70 after the 2-row insert) ----
$ws.Range("B70").Formula = "=SUM(A3,A5,A7,A9:A13,A15,A17:A20,A22:A25,A27:A28,A30:A31,A33:A39,A41:A43,A45:A47,A49,A51,A53,A55,A57,A59,A61,A63,A65,A67,A69)"

# --- Cosmetic: leave the view scrolled to where the new rows were typed --
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 1
[void]$ws.Range("B71").Select()

[void]$wb.Application.Calculate()
